$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.078393459320068
$ws.Range("B1").Value = 1.844270944595337
$ws.Range("C1").Value = 1.800463676452637
$ws.Range("D1").Value = 1.878517627716064
$ws.Range("E1").Value = 1.306295037269592
